$d = $word.ActiveDocument

# Locate the end of the paragraph that finishes with "...reconsiderări."
# (the last paragraph of the introductory red/italic guidance text,
# immediately before "Capitolul 1. Descrierea problemei").
$find = $d.Content
$found = $find.Find.Execute(
    "se vor sublinia mai ales aspecte ale cercetării în domeniu care nu au fost explorate sau care necesită investigații suplimentare sau reconsiderări.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the target paragraph (...reconsiderări.)"
}

$insertAt = $find.End

# Build a fresh Range object at that offset (reusing the Find range directly
# for InsertXML can clobber the paragraph it was found in, so re-anchor it).
$ins = $d.Range($insertAt, $insertAt)

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:color w:val="002060"/>
                <w:lang w:val="ro-RO" w:eastAsia="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="002060"/>
                <w:lang w:val="ro-RO" w:eastAsia="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">Odată cu trecerea timpului, tehnologia a jucat un rol esențial în ceea ce privește modul în care </w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

[void]$ins.InsertXML($xml)

Write-Host "Inserted new paragraph after the introductory guidance text."
